$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

$ws.Range("E2").Value = "RMA-HT43-001"
$ws.Range("F2").Value = "RMA-HT43-1-1"
$ws.Range("J2").Value = "a7s5f000000xLEkAAM"

$ws.Range("E3").Value = "RMA-HT43-002"
$ws.Range("F3").Value = "RMA-HT43-1-2"
$ws.Range("J3").Value = "a7s5f000000xLElAAM"

$ws.Range("E4").Value = "RMA-HT43-003"
$ws.Range("F4").Value = "RMA-HT43-1-3"
$ws.Range("J4").Value = "a7s5f000000xLEmAAM"
